$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$status = "Ready for handoff"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/oltest/blob/873fbfc7e6f717904da1091db5f762fa3f95778e/e2e/ad764da2-3354-4da9-831f-61edae4468bf.md, latest: https://github.com/OpenLocalizationTestOrg/oltest/blob/e51e6f0a667183ab00bbbb7cb6eccb9d85fd1508/e2e/ad764da2-3354-4da9-831f-61edae4468bf.md."

# --- Overview sheet: row 3 is the ad764da2-...md report row ---
# E3 = zh-cn status, F3 = de-de status, G3 = Latest HO Xliff Generate Date
$wsOverview.Range("E3").Value = $status
$wsOverview.Range("F3").Value = $status
$wsOverview.Range("G3").Value = "2016-08-13 07:00:38"

# --- zh-cn sheet: row 3 is the ad764da2-... file ---
# C3 = Status, H3 = Latest Handoff Datetime, P3 = Error Detail
$wsZhCn.Range("C3").Value = $status
$wsZhCn.Range("H3").Value = "2016-08-13 07:00:29"
$wsZhCn.Range("P3").Value = $errorDetail

# --- de-de sheet: row 3 is the ad764da2-... file ---
# C3 = Status, H3 = Latest Handoff Datetime, P3 = Error Detail
$wsDeDe.Range("C3").Value = $status
$wsDeDe.Range("H3").Value = "2016-08-13 07:00:38"
$wsDeDe.Range("P3").Value = $errorDetail

# Widen the Error Detail column (P) on both localized sheets to fit the long message
$wsZhCn.Columns.Item(16).ColumnWidth = 39.14
$wsDeDe.Columns.Item(16).ColumnWidth = 39.14
